# Bill of Materials update: swap the 8-position terminal block (J1) for a
# 10-position terminal block so the DUT can have 2 connections (one to
# connect the DUT, one to measure it), and refresh the affected supplier
# pricing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (J1 connector) ---------------------------------------------
# Name: CONN8 -> CONN10_2
$ws.Range("B2").Value = "CONN10_2"

# Description: 8 POS -> 10 POS
$ws.Range("C2").Value = "Conn Terminal Block 10 POS 2.54mm Solder ST Thru-Hole 6A"

# Manufacturer Part Number 1: OSTVN08A150 -> OSTVN10A150
$ws.Range("G2").Value = "OSTVN10A150"

# Supplier Part Number 1: ED10566-ND -> ED10567-ND
$ws.Range("J2").Value = "ED10567-ND"

# Supplier Unit Price 1 / Supplier Subtotal 1 (qty 1): 2.48 -> 3.1
$ws.Range("K2").Value = 3.1
$ws.Range("L2").Value = 3.1

# --- Row 4 (R1-R12 resistors) pricing refresh --------------------------
# Supplier Unit Price 1: 1.18 -> 1.31
$ws.Range("K4").Value = 1.31
# Supplier Subtotal 1 (qty 12): 14.21 -> 15.66
$ws.Range("L4").Value = 15.66
